$d = $word.ActiveDocument

# 1. Split "... pri pisanju uputstva za" into two runs:
#    "... pri pisanju uputstva " and "za "
$d.Content.Find.Execute(
    "uputstva za",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "uputstva ^pza ",
    2
)

# 2. Collapse the curly-quoted "Odjava" (three runs joined by proofErr
#    tags) into a single run "Odjava" within the smart quotes.
$d.Content.Find.Execute(
    [char]8220 + "Odjava" + [char]8221,
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    [char]8220 + "Odjava" + [char]8221,
    2
)
